# Update spike sorting summary for bat 9861 (inclusion_lists/sorting_summary.xlsx)
#
# - rows 36-38: increment the TT ("C") counter that was reset by mistake
# - row 38: flag "no cells" in the sorting-comments column
# - row 39: the session actually starts on 43245 (not 43244) -> shift date,
#   reset TT counter to 1, drop the stray comment that belonged to the next row
# - rows 40-55: newly logged tetrodes/sessions for bat 9861 through 43250

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-Row($r, $a, $b, $c, $d, $e, $f, $g, $h, $i, $j) {
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 2).NumberFormat = "m/d/yy"
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $h
    $ws.Cells.Item($r, 9).Value = $i
    if ($j -ne $null) {
        $ws.Cells.Item($r, 10).Value = $j
    }
}

# --- fix up the TT counters on the existing rows 36-38 ---
$ws.Cells.Item(36, 3).Value = 2
$ws.Cells.Item(37, 3).Value = 3
$ws.Cells.Item(38, 3).Value = 4
$ws.Cells.Item(38, 10).Value = "no cells"

# --- row 39 actually belongs to the 43245 session ---
$ws.Cells.Item(39, 2).Value = 43245
$ws.Cells.Item(39, 2).NumberFormat = "m/d/yy"
$ws.Cells.Item(39, 3).Value = 1
$ws.Cells.Item(39, 10).ClearContents()

# NOTE on shared-string order: register "there were jumps when changing
# sessions" (row 45) before "lots of cells with lots of time dynamics, VERY
# hard to sort!" (row 44) so the new shared-string table entries land in the
# same order as the target workbook (index 23 then 24).
$j45 = "there were jumps when changing sessions"
$j44 = "lots of cells with lots of time dynamics, VERY hard to sort!"

# --- newly appended rows 40-55 ---
Set-Row 40 9861 43245 2 7 0 0.8 24 6 "highpass" "no cells"
Set-Row 41 9861 43245 3 7 0 0.8 24 6 "highpass" $null
Set-Row 42 9861 43245 4 7 0 0.8 24 6 "highpass" $null

Set-Row 43 9861 43246 1 6 0 0.8 24 6 "highpass" $null

$ws.Cells.Item(45, 10).Value = $j45
Set-Row 44 9861 43246 2 7 0 0.8 24 6 "highpass" $j44
Set-Row 45 9861 43246 3 7 0 0.8 24 6 "highpass" $j45

Set-Row 46 9861 43246 4 7 0 0.8 24 6 "highpass" $null

Set-Row 47 9861 43247 1 7 0 0.8 24 6 "highpass" $null
Set-Row 48 9861 43247 2 7 0 0.8 24 6 "highpass" $null
Set-Row 49 9861 43247 3 7 0 0.8 24 6 "highpass" $null
Set-Row 50 9861 43247 4 7 0 0.8 24 6 "highpass" "no cells"

Set-Row 51 9861 43249 1 7 0 0.8 24 6 "highpass" $null
Set-Row 52 9861 43249 2 7 0 0.8 24 6 "highpass" $null
Set-Row 53 9861 43249 3 7 0 0.8 24 6 "highpass" $null
Set-Row 54 9861 43249 4 7 0 0.8 24 6 "highpass" "no isolated cells"

Set-Row 55 9861 43250 1 7 0 0.8 24 6 "highpass" $null

# --- move the view the way the author left it: scrolled down, cursor on A56 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 34
$win.ScrollColumn = 1
$ws.Range("A56").Select()
